$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1236520.6
$ws.Range("J76").Value = 5315.5
$ws.Range("L76").Value = 5315.5
$ws.Range("N76").Value = -5945.5
$ws.Range("H79").Value = 1236520.6
$ws.Range("J79").Value = 5315.5
$ws.Range("L79").Value = 5315.5
$ws.Range("N79").Value = -7499.5
$ws.Range("H98").Value = 2273.8096
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 1000
$ws.Range("N98").Value = -3996
$ws.Range("H112").Value = 4676.077
$ws.Range("J112").Value = 5032.4165
$ws.Range("L112").Value = 15097.2495
$ws.Range("N112").Value = -17313.2495
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").ClearContents()
$ws.Range("H122").Value = 2273.8096
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
$ws.Range("H129").Value = 918.36957
$ws.Range("J129").Value = 891.9318
$ws.Range("L129").Value = 2675.7954
$ws.Range("N129").Value = -12675.7954
$ws.Range("H132").Value = 1177.7742
$ws.Range("I132").Value = 1137.3334
$ws.Range("J132").Value = 1450.75
$ws.Range("K132").Value = 3412.0002
$ws.Range("L132").Value = 4352.25
$ws.Range("M132").Value = -882.0001999999999
$ws.Range("N132").Value = -9412.25
$ws.Range("H133").Value = 78899.5
$ws.Range("J133").Value = 78899.5
$ws.Range("L133").Value = 78899.5
$ws.Range("N133").Value = -89019.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2114761.2
$ws.Range("I2").Value = 3322762.5
$ws.Range("J2").Value = 759
$ws.Range("K2").Value = 3322762.5
$ws.Range("L2").Value = 759
$ws.Range("M2").Value = -3322649.5
$ws.Range("N2").Value = -985
$ws.Range("H61").Value = 5455.4614
$ws.Range("I61").Value = 3158
$ws.Range("K61").Value = 3158
$ws.Range("M61").Value = -2946
$ws.Range("H74").Value = 1223
$ws.Range("I74").Value = 783.5405
$ws.Range("K74").Value = 783.5405
$ws.Range("M74").Value = 90.45950000000005
$ws.Range("H77").Value = 1223
$ws.Range("I77").Value = 783.5405
$ws.Range("K77").Value = 3917.7025
$ws.Range("M77").Value = 450.2975000000001
$ws.Range("H116").Value = 2114761.2
$ws.Range("I116").Value = 3322762.5
$ws.Range("J116").Value = 759
$ws.Range("K116").Value = 3322762.5
$ws.Range("L116").Value = 759
$ws.Range("M116").Value = -3320468.5
$ws.Range("N116").Value = -5347
$ws.Range("H134").Value = 39532.668
$ws.Range("J134").Value = 39532.668
$ws.Range("L134").Value = 39532.668
$ws.Range("N134").Value = -49672.668
$ws.Range("H135").Value = 19000
$ws.Range("J135").Value = 19000
$ws.Range("L135").Value = 19000
$ws.Range("N135").Value = -29140
$ws.Range("H136").Value = 5455.4614
$ws.Range("I136").Value = 3158
$ws.Range("K136").Value = 9474
$ws.Range("M136").Value = -6924

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2114761.2
$ws.Range("I3").Value = 3322762.5
$ws.Range("J3").Value = 759
$ws.Range("K3").Value = 3322762.5
$ws.Range("L3").Value = 759
$ws.Range("M3").Value = -3322648.5
$ws.Range("N3").Value = -987
$ws.Range("H20").Value = 2174.8125
$ws.Range("I20").Value = 2196.7273
$ws.Range("K20").Value = 2196.7273
$ws.Range("M20").Value = -1949.7273
$ws.Range("H82").Value = 40875
$ws.Range("I82").Value = 20000
$ws.Range("K82").Value = 20000
$ws.Range("M82").Value = -19617
$ws.Range("H85").Value = 40875
$ws.Range("I85").Value = 20000
$ws.Range("K85").Value = 20000
$ws.Range("M85").Value = -18674
$ws.Range("H107").Value = 2875.875
$ws.Range("I107").Value = 2875.875
$ws.Range("K107").Value = 2875.875
$ws.Range("M107").Value = -955.875
$ws.Range("H131").Value = 47155.4
$ws.Range("J131").Value = 47155.4
$ws.Range("L131").Value = 47155.4
$ws.Range("N131").Value = -57235.4
$ws.Range("H135").Value = 27959.8
$ws.Range("J135").Value = 27959.8
$ws.Range("L135").Value = 27959.8
$ws.Range("N135").Value = -38099.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2718901
$ws.Range("I58").Value = 3953728
$ws.Range("K58").Value = 3953728
$ws.Range("M58").Value = -3953525
$ws.Range("H99").Value = 2241.3
$ws.Range("I99").Value = 1300
$ws.Range("J99").Value = 3182.6
$ws.Range("K99").Value = 1300
$ws.Range("L99").Value = 3182.6
$ws.Range("M99").Value = 198
$ws.Range("N99").Value = -6178.6
$ws.Range("H122").Value = 2046.3846
$ws.Range("I122").Value = 1991.1818
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 5973.5454
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -3523.5454
$ws.Range("N122").Value = -11950
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920
$ws.Range("H126").Value = 2241.3
$ws.Range("I126").Value = 1300
$ws.Range("J126").Value = 3182.6
$ws.Range("K126").Value = 3900
$ws.Range("L126").Value = 9547.799999999999
$ws.Range("M126").Value = -1430
$ws.Range("N126").Value = -14487.8
$ws.Range("H132").Value = 2190.0417
$ws.Range("I132").Value = 1171.875
$ws.Range("K132").Value = 3515.625
$ws.Range("M132").Value = -985.625
$ws.Range("H136").Value = 2718901
$ws.Range("I136").Value = 3953728
$ws.Range("K136").Value = 11861184
$ws.Range("M136").Value = -11858634

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 922.125
$ws.Range("J5").Value = 925.2857
$ws.Range("L5").Value = 2775.8571
$ws.Range("N5").Value = -2999.8571
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H105").Value = 3072.7273
$ws.Range("J105").Value = 3072.7273
$ws.Range("L105").Value = 9218.1819
$ws.Range("N105").Value = -14460.1819
$ws.Range("H122").Value = 1346.8334
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1316.2
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 11845.8
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -16745.8
$ws.Range("H135").Value = 922.125
$ws.Range("J135").Value = 925.2857
$ws.Range("L135").Value = 8327.5713
$ws.Range("N135").Value = -13397.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1864.9166
$ws.Range("I122").Value = 1553
$ws.Range("J122").Value = 2301.6
$ws.Range("K122").Value = 4659
$ws.Range("L122").Value = 6904.799999999999
$ws.Range("M122").Value = -2209
$ws.Range("N122").Value = -11804.8
$ws.Range("H126").Value = 1490655.2
$ws.Range("I126").Value = 1986588.6
$ws.Range("J126").Value = 102041.8
$ws.Range("K126").Value = 5959765.800000001
$ws.Range("L126").Value = 306125.4
$ws.Range("M126").Value = -5957295.800000001
$ws.Range("N126").Value = -311065.4
$ws.Range("H132").Value = 2408213
$ws.Range("I132").Value = 2961877.8
$ws.Range("K132").Value = 8885633.399999999
$ws.Range("M132").Value = -8883103.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3829.5
$ws.Range("I7").Value = 2787
$ws.Range("K7").Value = 2787
$ws.Range("M7").Value = -2675
$ws.Range("H40").Value = 7159.7
$ws.Range("I40").Value = 1933
$ws.Range("K40").Value = 1933
$ws.Range("M40").Value = -1797
$ws.Range("H46").Value = 1422.3914
$ws.Range("J46").Value = 1579.7368
$ws.Range("L46").Value = 1579.7368
$ws.Range("N46").Value = -1955.7368
$ws.Range("H55").Value = 320.09375
$ws.Range("I55").Value = 270.5
$ws.Range("J55").Value = 429.2
$ws.Range("K55").Value = 270.5
$ws.Range("L55").Value = 429.2
$ws.Range("M55").Value = -97.5
$ws.Range("N55").Value = -775.2
$ws.Range("H126").Value = 3829.5
$ws.Range("I126").Value = 2787
$ws.Range("K126").Value = 8361
$ws.Range("M126").Value = -5891

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2997.6667
$ws.Range("J62").Value = 2997.5
$ws.Range("L62").Value = 2997.5
$ws.Range("N62").Value = -4245.5
$ws.Range("H65").Value = 2997.6667
$ws.Range("J65").Value = 2997.5
$ws.Range("L65").Value = 14987.5
$ws.Range("N65").Value = -21227.5
$ws.Range("H100").Value = 578.7692
$ws.Range("I100").Value = 394.75
$ws.Range("K100").Value = 789.5
$ws.Range("M100").Value = -248.5
$ws.Range("H112").Value = 16776
$ws.Range("J112").Value = 16776
$ws.Range("L112").Value = 16776
$ws.Range("N112").Value = -19730
$ws.Range("H126").Value = 11511.357
$ws.Range("I126").Value = 16523.143
$ws.Range("J126").Value = 6499.5713
$ws.Range("K126").Value = 49569.429
$ws.Range("L126").Value = 19498.7139
$ws.Range("M126").Value = -47099.429
$ws.Range("N126").Value = -24438.7139

